$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data currently on the sheet.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Remove column A (the "Ticker Symbol" list) entirely; this shifts the
# "2016" earnings-per-share column (formerly column B) over to column A.
$ws.Columns.Item(1).Delete()

# Locate the row holding the smallest (most negative) earnings value among
# the remaining numeric data, now living in column A, rows 2..lastRow.
$minVal = $null
$minRow = $null
for ($r = 2; $r -le $lastRow; $r++) {
    $v = $ws.Cells.Item($r, 1).Value2
    if ($minVal -eq $null -or $v -lt $minVal) {
        $minVal = $v
        $minRow = $r
    }
}

# Highlight that cell with a solid red fill.
if ($minRow -ne $null) {
    $ws.Cells.Item($minRow, 1).Interior.Color = 255
}
